$d = $word.ActiveDocument

# Locate the last bullet item in the "Errores" list
# ("Falta agregar una propiedad ...") and remember its paragraph index.
$target = $null
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Falta agregar una propiedad que venga del back para trabar las opciones de respuesta del canto.*") {
        $target = $p
        $targetIndex = $i
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

# Insert a new paragraph right after it; Word clones the paragraph's
# formatting (pStyle + numPr list), matching the existing bullets.
$target.Range.InsertParagraphAfter()

# The newly inserted paragraph is now immediately after the target.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Con la mentira al final se pasa de puntos y no finaliza"
